$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new rows for Fortigate/Paloalto firewall configuration entries
$ws.Range("B16").Value = "FORTIGATE"
$ws.Range("B17").Value = "PALOALTO"
$ws.Range("C16").Value = "Fortigate 7.0"
$ws.Range("A16").Value = "Firewall"
$ws.Range("C17").Value = "palo"

$ws.Range("A17").Value = "Firewall"
$ws.Range("D16").Value = "CIS"
$ws.Range("D17").Value = "CIS"

$ws.Range("C17").Select()
